$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 13 de Mayo de 2020 a las 15:35
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 15:35"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1410168
$ws.Range("C4").Value = 1532
$ws.Range("D4").Value = 298593
$ws.Range("E4").Value = 1028084
$ws.Range("F4").Value = 16473
$ws.Range("G4").Value = 66
$ws.Range("H4").Value = 83491

# Row 7: Reino Unido
$ws.Range("A7").Value = "Reino Unido"
$ws.Range("B7").Value = 229705
$ws.Range("C7").Value = 3242
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 196175
$ws.Range("F7").Value = 1559
$ws.Range("G7").Value = 494
$ws.Range("H7").Value = 33186

# Row 11: Alemania
$ws.Range("A11").Value = "Alemania"
$ws.Range("B11").Value = 173546
$ws.Range("C11").Value = 375
$ws.Range("D11").Value = 148700
$ws.Range("E11").Value = 17066
$ws.Range("F11").Value = 1539
$ws.Range("G11").Value = 42
$ws.Range("H11").Value = 7780

# Row 72: Azerbaiyan
$ws.Range("A72").Value = "Azerbaiyan"
$ws.Range("B72").Value = 2758
$ws.Range("C72").Value = 65
$ws.Range("D72").Value = 1789
$ws.Range("E72").Value = 934
$ws.Range("F72").Value = 30
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 35

# Row 73: Grecia
$ws.Range("A73").Value = "Grecia"
$ws.Range("B73").Value = 2744
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 1374
$ws.Range("E73").Value = 1218
$ws.Range("F73").Value = 32
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 152

# Row 84: Islandia
$ws.Range("A84").Value = "Islandia"
$ws.Range("B84").Value = 1802
$ws.Range("C84").Value = 1
$ws.Range("D84").Value = 1780
$ws.Range("E84").Value = 12
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 10

# Row 111: Tayikistan
$ws.Range("A111").Value = "Tayikistan"
$ws.Range("B111").Value = 801
$ws.Range("C111").Value = 72
$ws.Range("D111").Value = 0
$ws.Range("E111").Value = 778
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 23

# Row 112: Burkina Faso
$ws.Range("A112").Value = "Burkina Faso"
$ws.Range("B112").Value = 766
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 588
$ws.Range("E112").Value = 127
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 51

# Row 113: Principado de Andorra
$ws.Range("A113").Value = "Principado de Andorra"
$ws.Range("B113").Value = 758
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 568
$ws.Range("E113").Value = 142
$ws.Range("F113").Value = 3
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 48

# Row 114: Paraguay
$ws.Range("A114").Value = "Paraguay"
$ws.Range("B114").Value = 737
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 173
$ws.Range("E114").Value = 554
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 10

# Row 115: Mali
$ws.Range("A115").Value = "Mali"
$ws.Range("B115").Value = 730
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 398
$ws.Range("E115").Value = 292
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 40

# Row 120: San Marino
$ws.Range("A120").Value = "San Marino"
$ws.Range("B120").Value = 643
$ws.Range("C120").Value = 5
$ws.Range("D120").Value = 161
$ws.Range("E120").Value = 441
$ws.Range("F120").Value = 2
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 41

# Row 122: Guinea Ecuatorial
$ws.Range("A122").Value = "Guinea Ecuatorial"
$ws.Range("B122").Value = 522
$ws.Range("C122").Value = 83
$ws.Range("D122").Value = 13
$ws.Range("E122").Value = 503
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 2
$ws.Range("H122").Value = 6

# Row 123: Tanzania
$ws.Range("A123").Value = "Tanzania"
$ws.Range("B123").Value = 509
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 183
$ws.Range("E123").Value = 305
$ws.Range("F123").Value = 7
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 21

# Row 124: Malta
$ws.Range("A124").Value = "Malta"
$ws.Range("B124").Value = 508
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 436
$ws.Range("E124").Value = 66
$ws.Range("F124").Value = 1
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 6

# Row 125: Jamaica
$ws.Range("A125").Value = "Jamaica"
$ws.Range("B125").Value = 507
$ws.Range("C125").Value = 2
$ws.Range("D125").Value = 100
$ws.Range("E125").Value = 398
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 9

# Row 126: Zambia
$ws.Range("A126").Value = "Zambia"
$ws.Range("B126").Value = 446
$ws.Range("C126").Value = 5
$ws.Range("D126").Value = 124
$ws.Range("E126").Value = 315
$ws.Range("F126").Value = 1
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 7

# Row 127: Taiwan
$ws.Range("A127").Value = "Taiwan"
$ws.Range("B127").Value = 440
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 375
$ws.Range("E127").Value = 58
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 7

# Row 177: Macao
$ws.Range("A177").Value = "Macao"
$ws.Range("B177").Value = 45
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 43
$ws.Range("E177").Value = 2
$ws.Range("F177").Value = 1
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 0

# Row 182: Zimbabue
$ws.Range("A182").Value = "Zimbabue"
$ws.Range("B182").Value = 37
$ws.Range("C182").Value = 1
$ws.Range("D182").Value = 12
$ws.Range("E182").Value = 21
$ws.Range("F182").Value = 0
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 4
